# A new daily price record is inserted into the "Kiwi" sheet, right
# after the existing row for A101 (market=11 / Vega Monumental
# Concepción / Bíobío). This pushes all the following rows (old
# 102..192) down by one (new 103..193), and the brand-new record lands
# in row 102.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 102, shifting rows 102-192 down to 103-193.
$ws.Rows.Item(102).Insert()

# Populate the newly inserted row 102 with the new data point.
$ws.Range("A102").Value = 11
$ws.Range("B102").Value = "Vega Monumental Concepción"
$ws.Range("C102").Value = "Bíobío"
$ws.Range("D102").Value = 44778
$ws.Range("E102").Value = 8
$ws.Range("F102").Value = "Fruta"
$ws.Range("G102").Value = 100101
$ws.Range("H102").Value = "Berries"
$ws.Range("I102").Value = 100101007
$ws.Range("J102").Value = "Kiwi"
$ws.Range("K102").Value = "Hayward"
$ws.Range("L102").Value = "Primera"
$ws.Range("M102").Value = 250
$ws.Range("N102").Value = 6000
$ws.Range("O102").Value = 6500
$ws.Range("P102").Value = 6300
$ws.Range("Q102").Value = "$/bandeja 18 kilos"
$ws.Range("R102").Value = "Provincia de Curicó"
$ws.Range("S102").Value = 350
$ws.Range("T102").Value = 18
